# Refresh the cryptos price table (GitHub Actions daily scrape update).
# Column D ("Price") values that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as literal text (matching the
# source data, which mixes "."-as-thousands-separator values like
# "29.476.22" with plain decimals like "242.00" - all stored as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.476.22"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").Value = "1.849.14"
$ws.Range("E3").Value = "  -0.62%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'242.00"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("D6").Value = "'0.6280"
$ws.Range("E6").Value = "  -2.10%  "

$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'48.11"
$ws.Range("E8").Value = "  +1.42%  "

$ws.Range("D9").Value = "'0.07544"

$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").Value = "'24.32"
$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.962.99"
$ws.Range("E12").Value = "  +5.52%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07697"
$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("D14").Value = "'5.005"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").Value = "'0.6861"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").Value = "'83.87"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").Value = "'0.000009746"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").Value = "2.165.62"
$ws.Range("E18").Value = "  +2.59%  "

$ws.Range("D19").Value = "'6.213"
$ws.Range("E19").Value = "  +1.53%  "

$ws.Range("D20").Value = "29.564.78"
$ws.Range("E20").Value = "  -0.81%  "

$ws.Range("D21").Value = "'234.45"
$ws.Range("E21").Value = "  -1.00%  "

$ws.Range("E22").Value = "  -1.54%  "

$ws.Range("D23").Value = "'0.9995"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'7.586"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "'155.31"
$ws.Range("E26").Value = "  -2.47%  "

$ws.Range("D27").Value = "'0.1393"
$ws.Range("E27").Value = "  -2.11%  "

$ws.Range("D28").Value = "'8.428"
$ws.Range("E28").Value = "  -1.63%  "

$ws.Range("E29").Value = "  -1.44%  "

$ws.Range("D30").Value = "'1.479"
$ws.Range("E30").Value = "  -1.46%  "

$ws.Range("D31").Value = "'0.05847"
$ws.Range("E31").Value = "  -5.43%  "

$ws.Range("E32").Value = "  -2.95%  "

$ws.Range("D33").Value = "'4.102"
$ws.Range("E33").Value = "  -1.50%  "

$ws.Range("D34").Value = "'4.019"
$ws.Range("E34").Value = "  -2.23%  "

$ws.Range("D35").Value = "'1.884"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("D36").Value = "'1.169"
$ws.Range("E36").Value = "  -0.69%  "

$ws.Range("D37").Value = "'0.7185"
$ws.Range("E37").Value = "  -1.80%  "

$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("D39").Value = "'2.793"
$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("D40").Value = "1.235.00"
$ws.Range("E40").Value = "  +1.88%  "

$ws.Range("D41").Value = "'0.01778"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").Value = "'0.9110"
$ws.Range("E42").Value = "  -1.45%  "

$ws.Range("D43").Value = "'6.124"
$ws.Range("E43").Value = "  -2.71%  "

$ws.Range("D44").Value = "2.074.58"
$ws.Range("E44").Value = "  +2.66%  "

$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").Value = "'101.79"
$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("E47").Value = "  +0.83%  "

$ws.Range("D48").Value = "'7.301"
$ws.Range("E48").Value = "  +8.68%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.162"
$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.4028"
$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.712"
$ws.Range("E51").Value = "  +2.21%  "
